$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Abbotsford"
$ws.Range("B2").Value = "Bodriggy Brewing Company  245 Johnston Street, Abbotsford VIC 3067"
$ws.Range("C2").Value = "28/12/20 2:50pm-5:30pm"
$ws.Range("D2").Value = "Case dined at venue"

$ws.Range("A3").Value = "Black Rock"
$ws.Range("B3").Value = "Smile Buffalo Thai restaurant  305 Beach Road, Black Rock VIC 3193"
$ws.Range("C3").Value = "27/12/20 7.30pm-9:00pm"
$ws.Range("D3").Value = "Case dined for dinner"

$ws.Range("A4").Value = "Brighton"
$ws.Range("B4").Value = "Royal Brighton Yacht Club - outdoor dining  253 Esplanade, Brighton VIC 3186"
$ws.Range("C4").Value = "29/12/20 12:00pm-2:00pm"
$ws.Range("D4").Value = "Case dined in outdoor restaurant"

$ws.Range("A5").Value = "Camberwell"
$ws.Range("B5").Value = "Tao Dumplings  1 Evans Place, Camberwell VIC 3124"
$ws.Range("C5").Value = "29/12/20 12:30pm-1:30pm"
$ws.Range("D5").Value = "Case ate at restaurant"

$ws.Range("A6").Value = "Collingwood"
$ws.Range("B6").Value = "Stomping Ground Brewing Company, 100 Gipps Street, Collingwood VIC 3066"
$ws.Range("C6").Value = "28/12/20 6:00pm-7:30pm"
$ws.Range("D6").Value = "Case dined at venue"

$ws.Range("A7").Value = "Docklands"
$ws.Range("B7").Value = "Melbourne Boat Hire - Yarra River Cruise. 45 Newquay Promenade, Docklands VIC 3008"
$ws.Range("C7").Value = "28/12/20 11.26am-2:00pm"
$ws.Range("D7").Value = "Case attended venue"

$ws.Range("A8").Value = "Doveton"
$ws.Range("B8").Value = "Holy Family Parish Doveton Catholic  100 Power Road, Doveton VIC 3177"
$ws.Range("C8").Value = "26/12/20 4:00pm-6:00pm"
$ws.Range("D8").Value = "Case attended Spanish Service"

$ws.Range("A9").Value = "Glen Waverley"
$ws.Range("B9").Value = "Village Century City  285-287 Springvale Road, Glen Waverley VIC 3150"
$ws.Range("C9").Value = "28/12/20 2:45pm-5:30pm"
$ws.Range("D9").Value = "Case attended Gold Class screening Wonder Woman 1984"

$ws.Range("A10").Value = "Hampton"
$ws.Range("B10").Value = "Merrymen Cafe, 2 Small Street, Hampton VIC"
$ws.Range("C10").Value = "28/12/20 1:30pm-2:30pm"
$ws.Range("D10").Value = "Case ate in store"

$ws.Range("A11").Value = "Keysborough"
$ws.Range("B11").Value = "Sikh Temple Keysborough, 198-206 Perry Road, Keysborough"
$ws.Range("C11").Value = "1/01/21 3:00pm-5:00pm"
$ws.Range("D11").Value = "Case visited venue"

$ws.Range("A12").Value = "McKinnon"
$ws.Range("B12").Value = "Hotlocks By Rachael Hairdresser, 260 McKinnon Road, McKinnon VIC 3204"
$ws.Range("C12").Value = "23/12/20 4:00pm-6:00pm"
$ws.Range("D12").Value = "Case had hair cut in store"

$ws.Range("A13").Value = "Melbourne"
$ws.Range("B13").Value = "Left Bank Melbourne, 1 Southbank Blvd"
$ws.Range("C13").Value = "25/12/20 12:00pm-3:00pm"
$ws.Range("D13").Value = "Case ate in store"

$ws.Range("A14").Value = "Melbourne"
$ws.Range("B14").Value = "Melbourne Central Lion Hotel, 211 La Trobe Street"
$ws.Range("C14").Value = "28/12/20 10:00pm-12.00am"
$ws.Range("D14").Value = "Case attended venue"

$ws.Range("A15").Value = "Melbourne"
$ws.Range("B15").Value = "Nandos  27 Elizabeth Street, Melbourne"
$ws.Range("C15").Value = "01/01/2021 2:00am-2:30am"
$ws.Range("D15").Value = "Case dined at venue"

$ws.Range("A16").Value = "Moorabbin"
$ws.Range("B16").Value = "Grape and Grain Liquor Cellars, 14/16 Station St"
$ws.Range("C16").Value = "24/12/20 1:00pm-10:00pm  28/12/20 8.05pm-8.47pm  29/12/20 12:00pm-4:00pm"
$ws.Range("D16").Value = "Case's workplace"

$ws.Range("A17").Value = "Mordialloc"
$ws.Range("B17").Value = "Woodlands Golf Club - club bar  109 White Street Mordialloc VIC 3195"
$ws.Range("C17").Value = "23/12/20 12:30pm-1:30pm"
$ws.Range("D17").Value = "Case attended club house bar"

$ws.Range("A18").Value = "Mordialloc"
$ws.Range("B18").Value = "Woodlands Golf Club - club bar  109 White Street Mordialloc VIC 3195"
$ws.Range("C18").Value = "28/12/20 4:40pm-5:15pm"
$ws.Range("D18").Value = "Case attended club house bar"

$ws.Range("A19").Value = "South Melbourne"
$ws.Range("B19").Value = "The Nike Company  134 Buckhurst Street, South Melbourne"
$ws.Range("C19").Value = "30/12/20, 12:00pm-12:45pm"
$ws.Range("D19").Value = "Case shopped"

$ws.Range("A20").Value = "Southbank"
$ws.Range("B20").Value = "Rockpool Bar and Grill, Crown Casino  8 Whiteman Street, Southbank"
$ws.Range("C20").Value = "23/12/20 8:00pm-11:00pm"
$ws.Range("D20").Value = "Case ate in store"

$ws.Range("A21").Value = "Springvale"
$ws.Range("B21").Value = "IKEA Springvale - Cafe and Restaurant, 917 Princes Hwy"
$ws.Range("C21").Value = "30/12/20 5:30pm-6.30pm"
$ws.Range("D21").Value = "Case dined at cafe"
